# sr: update site form
#
# - Appends 6 new "community_id" choice rows (T1 Spare 4-6, T2 Spare 4-6,
#   ids 119-124) right after the existing list on the "choices" sheet.
# - Leaves the "choices" sheet scrolled to/selected near the newly added
#   rows, and makes "choices" the active sheet/tab of the workbook (the
#   "survey" sheet was active before).

$wb = $excel.ActiveWorkbook
$choices = $wb.Worksheets.Item("choices")

$newRows = @(
    @{ Id = 119; Label = "T1 Spare 4" },
    @{ Id = 120; Label = "T1 Spare 5" },
    @{ Id = 121; Label = "T1 Spare 6" },
    @{ Id = 122; Label = "T2 Spare 4" },
    @{ Id = 123; Label = "T2 Spare 5" },
    @{ Id = 124; Label = "T2 Spare 6" }
)

$startRow = 70
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $choices.Cells.Item($r, 1).Value = "community_id"
    $choices.Cells.Item($r, 2).Value = $row.Id
    $choices.Cells.Item($r, 3).Value = $row.Id
    # Column D stays blank for these rows (matches the blank-but-present
    # D cells used throughout the rest of the community_id block).
    $choices.Cells.Item($r, 4).Style = "Normal"
    $choices.Cells.Item($r, 5).Value = $row.Label
}

# Make "choices" the active sheet and select/scroll near the new rows,
# reflecting the state the workbook was left in after the edit.
$choices.Activate()
$choices.Range("A76:XFD76").Select() | Out-Null
